# Adding changes for MTTR Suite
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the sample data populated in A2:E4 (the MTTR suite no longer ships
# canned record ids in this datasheet). Use Clear() so the cells are dropped
# entirely (value + style) rather than left behind as blank styled cells --
# except D2, which keeps its distinct style ("s"=2) and is only emptied.
$ws.Range("A2:C2").Clear()
$ws.Range("D2").ClearContents()
$ws.Range("E2").Clear()
$ws.Range("A3:E3").Clear()
$ws.Range("A4").Clear()

# Update the active selection to match the new state.
$ws.Range("A2:E5").Select()
